$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of (Attribute, Type) pairs for rows 2..21.
# The set of pairs is unchanged; only the row order changes.
$rows = @(
    @("operation_end_time", "datetime"),
    @("parameters", "dict"),
    @("SubProcessID", "str"),
    @("identifier:id", "str"),
    @("lifecycle:transition", "str"),
    @("complete_service_time", "str"),
    @("concept:name", "str"),
    @("process_model_id", "str"),
    @("event_id", "str"),
    @("case", "str"),
    @("org:resource", "str"),
    @("requested_service_url", "str"),
    @("response_status_code", "float"),
    @("unsatisfied_condition_description", "str"),
    @("human_workstation_green_button_pressed", "float"),
    @("time:timestamp", "datetime"),
    @("planned_operation_time", "str"),
    @("lifecycle:state", "str"),
    @("case:concept:name", "str"),
    @("current_task", "str")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}
